$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'62.807.36"
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").Value = "'2.941.40"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'593.30"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").Value = "'147.15"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.76%  "

$ws.Range("D9").Value = "'2.939.60"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "'7.32"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("E11").Value = "  +6.09%  "

$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "'0.0000238"
$ws.Range("E13").Value = "  +5.51%  "

$ws.Range("D14").Value = "'32.79"
$ws.Range("E14").Value = "  -2.72%  "

$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "'3.425.43"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "'62.746.18"
$ws.Range("E17").Value = "  +2.13%  "

$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'2.960.66"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").Value = "'441.09"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D21").Value = "'13.41"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("D22").Value = "'0.666"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").Value = "'7.04"
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").Value = "'81.30"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").Value = "'11.12"
$ws.Range("E25").Value = "  +1.70%  "

$ws.Range("D26").Value = "'2.14"
$ws.Range("E26").Value = "  -3.02%  "

$ws.Range("D27").Value = "'11.75"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = "  +4.22%  "

$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("D32").Value = "'0.0000103"
$ws.Range("E32").Value = "  +15.68%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.109"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.39"
$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = "  -2.11%  "

$ws.Range("D37").Value = "'3.13"
$ws.Range("E37").Value = "  +4.13%  "

$ws.Range("E38").Value = "  -1.17%  "

$ws.Range("D39").Value = "'49.65"
$ws.Range("E39").Value = "  -0.68%  "

$ws.Range("E40").Value = "  +0.94%  "

$ws.Range("D41").Value = "'8.49"
$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("E42").Value = "  -5.52%  "

$ws.Range("D43").Value = "'0.281"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").Value = "'39.44"
$ws.Range("E44").Value = "  -7.13%  "

$ws.Range("D45").Value = "'2.702.93"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'134.59"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "'362.75"
$ws.Range("E47").Value = "  -0.65%  "

$ws.Range("E48").Value = "  -3.40%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").Value = "'22.81"
$ws.Range("E51").Value = "  -4.25%  "
